# Auto-generated Excel COM-interop script to apply F-column ("想去人数") updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(3, 6).Value = 1322   # F3: 1318 -> 1322
$ws.Cells.Item(4, 6).Value = 1119   # F4: 1116 -> 1119
$ws.Cells.Item(5, 6).Value = 1007   # F5: 1005 -> 1007
$ws.Cells.Item(6, 6).Value = 1784   # F6: 1781 -> 1784
$ws.Cells.Item(7, 6).Value = 551   # F7: 550 -> 551
$ws.Cells.Item(8, 6).Value = 1187   # F8: 1183 -> 1187
$ws.Cells.Item(9, 6).Value = 55   # F9: 54 -> 55
$ws.Cells.Item(11, 6).Value = 125   # F11: 123 -> 125
$ws.Cells.Item(12, 6).Value = 286   # F12: 283 -> 286
$ws.Cells.Item(13, 6).Value = 63   # F13: 62 -> 63
$ws.Cells.Item(14, 6).Value = 89   # F14: 86 -> 89
$ws.Cells.Item(15, 6).Value = 677   # F15: 670 -> 677
$ws.Cells.Item(16, 6).Value = 163   # F16: 159 -> 163
$ws.Cells.Item(21, 6).Value = 139   # F21: 138 -> 139
$ws.Cells.Item(22, 6).Value = 664   # F22: 662 -> 664
$ws.Cells.Item(23, 6).Value = 30   # F23: 28 -> 30
$ws.Cells.Item(25, 6).Value = 147   # F25: 146 -> 147
$ws.Cells.Item(28, 6).Value = 311   # F28: 310 -> 311
$ws.Cells.Item(29, 6).Value = 156   # F29: 154 -> 156
$ws.Cells.Item(31, 6).Value = 268   # F31: 266 -> 268

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(4, 6).Value = 316   # F4: 315 -> 316
$ws.Cells.Item(7, 6).Value = 251   # F7: 250 -> 251
$ws.Cells.Item(8, 6).Value = 82   # F8: 81 -> 82

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 307   # F2: 306 -> 307

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 307   # F2: 306 -> 307
$ws.Cells.Item(4, 6).Value = 1322   # F4: 1318 -> 1322
$ws.Cells.Item(5, 6).Value = 1119   # F5: 1116 -> 1119
$ws.Cells.Item(6, 6).Value = 1007   # F6: 1005 -> 1007
$ws.Cells.Item(7, 6).Value = 1784   # F7: 1781 -> 1784
$ws.Cells.Item(8, 6).Value = 552   # F8: 550 -> 552
$ws.Cells.Item(9, 6).Value = 1187   # F9: 1183 -> 1187
$ws.Cells.Item(10, 6).Value = 55   # F10: 54 -> 55
$ws.Cells.Item(13, 6).Value = 125   # F13: 123 -> 125
$ws.Cells.Item(14, 6).Value = 286   # F14: 283 -> 286
$ws.Cells.Item(15, 6).Value = 63   # F15: 62 -> 63
$ws.Cells.Item(16, 6).Value = 89   # F16: 86 -> 89
$ws.Cells.Item(17, 6).Value = 677   # F17: 670 -> 677
$ws.Cells.Item(18, 6).Value = 163   # F18: 159 -> 163
$ws.Cells.Item(22, 6).Value = 316   # F22: 315 -> 316
$ws.Cells.Item(27, 6).Value = 251   # F27: 250 -> 251
$ws.Cells.Item(28, 6).Value = 251   # F28: 250 -> 251
$ws.Cells.Item(29, 6).Value = 139   # F29: 138 -> 139
$ws.Cells.Item(30, 6).Value = 664   # F30: 662 -> 664
$ws.Cells.Item(31, 6).Value = 30   # F31: 28 -> 30
$ws.Cells.Item(33, 6).Value = 147   # F33: 146 -> 147
$ws.Cells.Item(36, 6).Value = 311   # F36: 310 -> 311
$ws.Cells.Item(37, 6).Value = 82   # F37: 81 -> 82
$ws.Cells.Item(39, 6).Value = 156   # F39: 154 -> 156
$ws.Cells.Item(41, 6).Value = 268   # F41: 266 -> 268
